# Automatische test-sync: 2025-07-22 12:14:50
# Adds the second test-mail row ("Bestelling / Levering") to the Logs
# sheet, mirrors the per-category tally on the Dashboard sheet, and
# extends the conditional-formatting ranges + chart series references
# so they cover the newly added row.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------
# 1) Logs sheet: append row 3 with the new test mail
# ---------------------------------------------------------------
$logs = $wb.Worksheets.Item("Logs")

$logs.Range("A3").Value = "Kun je deze order vandaag nog verwerken?"
$logs.Range("B3").Value = "mailmind.test@zohomail.eu"
$logs.Range("C3").Value = "Testmail #2: Kun je deze order vandaag nog verwerken?"
$logs.Range("D3").Value = "Bestelling / Levering"
$logs.Range("E3").Value = "Beste klant,`nDank u voor uw bericht. Om uw vraag over het verwerken van de order vandaag te beantwoorden, hebben we wat meer informatie nodig. Kunt u ons het ordernummer doorgeven, zodat we kunnen nakijken of het mogelijk is om de order vandaag nog te verwerken?`nMet vriendelijke groet,`n[Bedrijfsnaam] E-mailassistent"
$logs.Range("F3").Value = "2025-07-22 12:13:52"
$logs.Range("G3").Value = "Ja"
$logs.Range("H3").Value = "Ja"
$logs.Range("I3").Value = "Nee"
$logs.Range("J3").Value = "Nee"

# Writing the multi-line reply text auto-marks row 3 with a custom
# height; AutoFit() recomputes it as a normal (non-custom) row height,
# same as row 2, so the row picks up the sheet's default sizing again.
$logs.Rows.Item(3).AutoFit()

# Extend the conditional-formatting blocks (one per column) so they
# also cover row 3, same as the existing row-2 rule.
foreach ($col in @("D", "G", "H", "I", "J")) {
    $oldRange = $logs.Range($col + "2")
    $newRange = $logs.Range($col + "2:" + $col + "3")
    $fcs = $oldRange.FormatConditions
    for ($i = 1; $i -le $fcs.Count(); $i++) {
        $fcs.Item($i).ModifyAppliesToRange($newRange)
    }
}

# ---------------------------------------------------------------
# 2) Dashboard sheet: append row 3 with the updated tally
# ---------------------------------------------------------------
$dash = $wb.Worksheets.Item("Dashboard")

$dash.Range("A3").Value = "Bestelling / Levering"
$dash.Range("B3").Value = 1

# ---------------------------------------------------------------
# 3) Chart: widen the category/value series references to include
#    the new Dashboard row.
# ---------------------------------------------------------------
$chartObj = $dash.ChartObjects().Item(1)
$chart = $chartObj.Chart
$series = $chart.SeriesCollection().Item(1)
$series.Formula = "=SERIES(Dashboard!`$B`$1,Dashboard!`$A`$2:`$A`$3,Dashboard!`$B`$2:`$B`$3,1)"
